# Re-sort the bibliography by Type (column E) then Year (column D),
# tighten up the long text columns, let row heights auto-fit again,
# and refresh the view/zoom to where I was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A1:H36")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("E2:E36"), 0, 1, 0, 0) | Out-Null
$ws.Sort.SortFields.Add($ws.Range("D2:D36"), 0, 1, 0, 0) | Out-Null

$ws.Sort.SetRange($dataRange)
$ws.Sort.Header = 1
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = 1
$ws.Sort.Apply()

# Narrower columns now that entries have been tidied up
$ws.Columns.Item(5).ColumnWidth = 8.42578125
$ws.Columns.Item(6).ColumnWidth = 35.42578125
$ws.Columns.Item(7).ColumnWidth = 40.7109375
$ws.Columns.Item(8).ColumnWidth = 25.140625

# Let the rows auto-fit again after the column/sort changes
$ws.Rows.Item("2:36").EntireRow.AutoFit()

# Restore the view: scrolled down a bit, zoomed out slightly
$excel.ActiveWindow.Zoom = 85
$ws.Range("F27").Select()
$excel.ActiveWindow.ScrollRow = 26
